$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("ECs", "Bsg", "Sele", "ECs", 3, 1, 44.09582266666666, 132.287468, 0.1927468402671175, 0.1927468402671175, 3, 1, 4.699506666666667, 14.09852, 0.9660495246229048, 0.9660495246229047, 207.2286125941511, 1865.05751334736, 0.1862029934126159, 0.1862029934126158)
    ,@("ECs", "Bsg", "Sele", "FAPs", 3, 1, 44.09582266666666, 132.287468, 0.1927468402671175, 0.1927468402671175, 1, 0.3333333333333333, 0.1651576666666667, 0.495473, 0.03395047537709522, 0.03395047537709522, 7.282763181373777, 65.544868632364, 0.006543846854501678, 0.006543846854501678)
    ,@("FAPs", "Bsg", "Sele", "ECs", 3, 1, 111.3149186666667, 333.944756, 0.4865676055026886, 0.4865676055026886, 3, 1, 4.699506666666667, 14.09852, 0.9660495246229048, 0.9660495246229047, 523.1252023734578, 4708.12682136112, 0.4700484039927774, 0.4700484039927773)
    ,@("FAPs", "Bsg", "Sele", "FAPs", 3, 1, 111.3149186666667, 333.944756, 0.4865676055026886, 0.4865676055026886, 1, 0.3333333333333333, 0.1651576666666667, 0.495473, 0.03395047537709522, 0.03395047537709522, 18.38451223217644, 165.460610089588, 0.01651920150991121, 0.01651920150991121)
    ,@("M2", "Bsg", "Sele", "ECs", 3, 1, 40.34450033333334, 121.033501, 0.176349470111689, 0.176349470111689, 3, 1, 4.699506666666667, 14.09852, 0.9660495246229048, 0.9660495246229047, 189.5992482798356, 1706.39323451852, 0.1703623217688983, 0.1703623217688983)
    ,@("M2", "Bsg", "Sele", "FAPs", 3, 1, 40.34450033333334, 121.033501, 0.176349470111689, 0.176349470111689, 1, 0.3333333333333333, 0.1651576666666667, 0.495473, 0.03395047537709522, 0.03395047537709522, 6.66320353788589, 59.968831840973, 0.005987148342790685, 0.005987148342790685)
    ,@("sCs", "Bsg", "Sele", "ECs", 3, 1, 33.02061066666666, 99.061832, 0.1443360841185049, 0.144336084118505, 3, 1, 4.699506666666667, 14.09852, 0.9660495246229048, 0.9660495246229047, 155.1805799654044, 1396.62521968864, 0.1394358054486133, 0.1394358054486133)
    ,@("sCs", "Bsg", "Sele", "FAPs", 3, 1, 33.02061066666666, 99.061832, 0.1443360841185049, 0.144336084118505, 1, 0.3333333333333333, 0.1651576666666667, 0.495473, 0.03395047537709522, 0.03395047537709522, 5.45360700961511, 49.08246308653599, 0.004900278669891646, 0.004900278669891647)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $col = 1 + $c
        $ws.Cells.Item($row, $col).Value = $rowData[$c]
    }
}

Write-Host "done"